$wb = $excel.ActiveWorkbook

# --- Update the Main sheet: Price (C3) changes from 101.45 to 100.70 ---
# Dependent formulas (Mkt Cap in C5, EV in C8) will recalculate automatically.
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("C3").Value = 100.7
$wsMain.Range("C4").Select() | Out-Null

# --- Rename the "DCF" sheet to "Model" ---
$wsModel = $wb.Worksheets.Item("DCF")
$wsModel.Name = "Model"

# --- Add new expense/gain line items to the Model sheet ---
$wsModel.Range("B16").Value = "Gain on Extinguishment"
$wsModel.Range("B17").Value = "Environmental Expense"
$wsModel.Range("B18").Value = "Pension Plan Termination Expense"

# Re-fit column B so the new, longer label is fully visible
$wsModel.Columns("B").AutoFit() | Out-Null

# Update page setup (paper size / orientation) on the Model sheet
$wsModel.PageSetup.PaperSize = 9
$wsModel.PageSetup.Orientation = 1

# Activate Model sheet and update its selection
$wsModel.Activate() | Out-Null
$wsModel.Range("B19").Select() | Out-Null

# --- Remove the now-unused "DDM" sheet ---
$wsDdm = $wb.Worksheets.Item("DDM")
$wsDdm.Delete() | Out-Null
